$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header cells (H1:P1) with the new column names ---
$headers = @("grade_total","grade_distance","grade_visitation","grade_encounters","NEVER","RARELY","SOMETIMES","FREQUENTLY","ALWAYS")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 8 + $i).Value = $headers[$i]
}

# --- Values shared by every data row (H:P), in column order ---
$rowValues = @(0, 0, 0, 0, 1.032, 1.023, 1.06, 1.128, 1.756)

# --- Swap columns A (originally year) and B (originally month) so ---
# --- that A becomes month and B becomes year, and populate the new ---
# --- H:P columns for every data row                                 ---
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $origAVal = $ws.Cells.Item($r, 1).Value2   # originally year
    $origBVal = $ws.Cells.Item($r, 2).Value2   # originally month

    $ws.Cells.Item($r, 1).Value = $origBVal    # A -> month
    $ws.Cells.Item($r, 2).Value = $origAVal    # B -> year

    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item($r, 8 + $i).Value = $rowValues[$i]
    }
}
